$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (A21) holding the same text value as A4/A... ("1234567890").
# Build it as a formula that evaluates to a text string, then convert the
# formula to a static value via copy/paste-special so the cell ends up a
# plain shared-string text cell (no numeric coercion, no style changes).
$ws.Range("A21").Formula = '="1234567890"'
$ws.Range("A21").Copy()
$ws.Range("A21").PasteSpecial(-4163)
$excel.CutCopyMode = 0
